$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8046
$ws1.Range("F5").Value = 5863
$ws1.Range("F6").Value = 498
$ws1.Range("F7").Value = 88
$ws1.Range("F10").Value = 289
$ws1.Range("F11").Value = 395
$ws1.Range("F12").Value = 66

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8046
$ws4.Range("F5").Value = 5863
$ws4.Range("F6").Value = 498
$ws4.Range("F7").Value = 88
$ws4.Range("F10").Value = 289
$ws4.Range("F14").Value = 395
$ws4.Range("F15").Value = 66
